# Agrego explicacion de como ejecutar el programa
#
# 1) Add a new paragraph (right after the last math formula, right before
#    the "Ingreso de Datos" heading) explaining how to run the app.
# 2) Move the "_GoBack" bookmark from the end of the "Visualizar tabla de
#    sumatorias" intro paragraph into the middle of the word "pares" in the
#    "Ingreso de Datos" section (this is what Word does automatically while
#    editing - the bookmark marks the last edited spot).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: insert the new explanatory paragraph.
# ---------------------------------------------------------------------
$headingRange = $d.Content
$headingRange.Find.Execute("Ingreso de Datos") | Out-Null
$headingPara = $headingRange.Paragraphs(1)
$emptyPara = $headingPara.Previous()

$ip = $emptyPara.Range
$ip.Collapse(1)                                   # wdCollapseStart

$ip.InsertAfter("Para iniciar la aplicación se debe ejecutar el script ")
$ip.Collapse(0)                                   # wdCollapseEnd

$ip.InsertAfter("UI.m")
$ip.Collapse(0)

$ip.InsertAfter(" con Octave.")

# ---------------------------------------------------------------------
# Step 2: relocate the "_GoBack" bookmark.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$splitRange = $d.Content
$splitRange.Find.Execute("cantidad de pares de números se desee ingresar.") | Out-Null
$splitPos = $splitRange.Start + "cantidad de p".Length
$newBookmarkRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
